$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp shown in F1
$ws.Range("F1").Value = "Last status check on: 02.02.2022 09:02"

# Update row 8 (Benzina Albert Modrice) prices
$ws.Range("B8").Value = 36.5
$ws.Range("C8").Value = 36.4

# D8 and E8 switch from numeric (price delta / serial date) to plain text values.
# Force a text number format before assigning so Excel does not reinterpret
# "+0.1" as a number or the date/time string as a date serial, then clear the
# explicit style back to the default (Normal) so no custom numFmt lingers.
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "+0.1"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2022-02-02 09:14:07"
$ws.Range("E8").Style = "Normal"
